$d = $word.ActiveDocument
$d.Content.Find.Execute("Chi phí phát triển  +  Chi phí kiểm thử: 20.000.000 VND", $true, $false, $false, $false, $false, $true, 1, $false, "Chi phí phát triển  +  Chi phí kiểm thử: 12.000.000 VND", 2)
